$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.153049230575562
$ws.Range("B1").Value = 2.576439142227173
$ws.Range("C1").Value = 5.978178977966309
$ws.Range("D1").Value = 2.125722646713257
$ws.Range("E1").Value = 1.223624587059021
